# Auto-generated market-data refresh for FFXIV Yojimbo profit sheets.
# Mirrors a scheduled-runner update to columns H:N (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 760
$ws.Range("I70").Value = 600
$ws.Range("J70").Value = 789.0909
$ws.Range("K70").Value = 1800
$ws.Range("L70").Value = 2367.2727
$ws.Range("M70").Value = -1530
$ws.Range("N70").Value = -2907.2727
$ws.Range("H73").Value = 760
$ws.Range("I73").Value = 600
$ws.Range("J73").Value = 789.0909
$ws.Range("K73").Value = 1800
$ws.Range("L73").Value = 2367.2727
$ws.Range("M73").Value = -864
$ws.Range("N73").Value = -4239.2727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3250.2
$ws.Range("I61").Value = 1655
$ws.Range("J61").Value = 5199.8887
$ws.Range("K61").Value = 1655
$ws.Range("L61").Value = 5199.8887
$ws.Range("M61").Value = -1443
$ws.Range("N61").Value = -5623.8887
$ws.Range("H74").Value = 1930.2222
$ws.Range("I74").Value = 1011.0769
$ws.Range("J74").Value = 4320
$ws.Range("K74").Value = 1011.0769
$ws.Range("L74").Value = 4320
$ws.Range("M74").Value = -137.0769
$ws.Range("N74").Value = -6068
$ws.Range("H77").Value = 1930.2222
$ws.Range("I77").Value = 1011.0769
$ws.Range("J77").Value = 4320
$ws.Range("K77").Value = 5055.3845
$ws.Range("L77").Value = 21600
$ws.Range("M77").Value = -687.3845000000001
$ws.Range("N77").Value = -30336
$ws.Range("H88").Value = 2793.75
$ws.Range("J88").Value = 5000
$ws.Range("L88").Value = 5000
$ws.Range("N88").Value = -5812
$ws.Range("H91").Value = 2793.75
$ws.Range("J91").Value = 5000
$ws.Range("L91").Value = 5000
$ws.Range("N91").Value = -7808
$ws.Range("H132").Value = 2359.7463
$ws.Range("I132").Value = 2068.6726
$ws.Range("J132").Value = 3693.8333
$ws.Range("K132").Value = 6206.0178
$ws.Range("L132").Value = 11081.4999
$ws.Range("M132").Value = -3676.0178
$ws.Range("N132").Value = -16141.4999
$ws.Range("H136").Value = 3250.2
$ws.Range("I136").Value = 1655
$ws.Range("J136").Value = 5199.8887
$ws.Range("K136").Value = 4965
$ws.Range("L136").Value = 15599.6661
$ws.Range("M136").Value = -2415
$ws.Range("N136").Value = -20699.6661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2451
$ws.Range("I20").Value = 2616.5
$ws.Range("J20").Value = 2351.7
$ws.Range("K20").Value = 2616.5
$ws.Range("L20").Value = 2351.7
$ws.Range("M20").Value = -2369.5
$ws.Range("N20").Value = -2845.7
$ws.Range("H48").Value = 104900
$ws.Range("J48").Value = 104900
$ws.Range("L48").Value = 104900
$ws.Range("N48").Value = -105730
$ws.Range("H86").Value = 4019.9592
$ws.Range("I86").Value = 3604.0488
$ws.Range("J86").Value = 6151.5
$ws.Range("K86").Value = 3604.0488
$ws.Range("L86").Value = 6151.5
$ws.Range("M86").Value = -2481.0488
$ws.Range("N86").Value = -8397.5
$ws.Range("H89").Value = 4019.9592
$ws.Range("I89").Value = 3604.0488
$ws.Range("J89").Value = 6151.5
$ws.Range("K89").Value = 18020.244
$ws.Range("L89").Value = 30757.5
$ws.Range("M89").Value = -12404.244
$ws.Range("N89").Value = -41989.5
$ws.Range("H134").Value = 1736.125
$ws.Range("I134").Value = 1570.4546
$ws.Range("K134").Value = 4711.3638
$ws.Range("M134").Value = -2176.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 968.5
$ws.Range("I16").Value = 962.2
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 962.2
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -675.2
$ws.Range("N16").Value = -1574
$ws.Range("H58").Value = 2446.818
$ws.Range("I58").Value = 3269.5833
$ws.Range("J58").Value = 1459.5
$ws.Range("K58").Value = 3269.5833
$ws.Range("L58").Value = 1459.5
$ws.Range("M58").Value = -3066.5833
$ws.Range("N58").Value = -1865.5
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -18740
$ws.Range("H113").Value = 968.5
$ws.Range("I113").Value = 962.2
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 962.2
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1207.8
$ws.Range("N113").Value = -5340
$ws.Range("H136").Value = 2446.818
$ws.Range("I136").Value = 3269.5833
$ws.Range("J136").Value = 1459.5
$ws.Range("K136").Value = 9808.749899999999
$ws.Range("L136").Value = 4378.5
$ws.Range("M136").Value = -7258.749899999999
$ws.Range("N136").Value = -9478.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 140.06667
$ws.Range("I40").Value = 140.06667
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 560.26668
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -491.26668
$ws.Range("N40").Value = $null
$ws.Range("H42").Value = 2155.3845
$ws.Range("I42").Value = 300
$ws.Range("J42").Value = 2310
$ws.Range("K42").Value = 900
$ws.Range("L42").Value = 6930
$ws.Range("M42").Value = -366
$ws.Range("N42").Value = -7998
$ws.Range("H80").Value = 3791.4546
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 4022.889
$ws.Range("K80").Value = 8250
$ws.Range("L80").Value = 12068.667
$ws.Range("M80").Value = -7314
$ws.Range("N80").Value = -13940.667
$ws.Range("H83").Value = 3791.4546
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 4022.889
$ws.Range("K83").Value = 24750
$ws.Range("L83").Value = 36206.001
$ws.Range("M83").Value = -20070
$ws.Range("N83").Value = -45566.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4067.6667
$ws.Range("I11").Value = 1751.5
$ws.Range("J11").Value = 8700
$ws.Range("K11").Value = 1751.5
$ws.Range("L11").Value = 8700
$ws.Range("M11").Value = -1612.5
$ws.Range("N11").Value = -8978
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352
$ws.Range("H132").Value = 1653.7307
$ws.Range("I132").Value = 1255.2307
$ws.Range("J132").Value = 2849.2307
$ws.Range("K132").Value = 3765.6921
$ws.Range("L132").Value = 8547.6921
$ws.Range("M132").Value = -1235.6921
$ws.Range("N132").Value = -13607.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4217.1113
$ws.Range("I62").Value = 3590.8
$ws.Range("K62").Value = 3590.8
$ws.Range("M62").Value = -2966.8
$ws.Range("H65").Value = 4217.1113
$ws.Range("I65").Value = 3590.8
$ws.Range("K65").Value = 17954
$ws.Range("M65").Value = -14834
$ws.Range("H122").Value = 627029.9399999999
$ws.Range("I122").Value = 2001199.8
$ws.Range("J122").Value = 2407.2727
$ws.Range("K122").Value = 6003599.4
$ws.Range("L122").Value = 7221.8181
$ws.Range("M122").Value = -6001149.4
$ws.Range("N122").Value = -12121.8181
$ws.Range("H136").Value = 218392.17
$ws.Range("I136").Value = 243361.47
$ws.Range("J136").Value = 1991.6666
$ws.Range("K136").Value = 730084.41
$ws.Range("L136").Value = 5974.9998
$ws.Range("M136").Value = -727534.41
$ws.Range("N136").Value = -11074.9998
